# Automatische test-sync: 2025-08-14 21:58:50
# Appends a new "Logs" row (row 34) describing an internal order request,
# extends the conditional-formatting ranges to cover it, and bumps the
# "Intern verzoek / Actie voor medewerker" tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 34

$logs.Range("A$newRow").Value = "Nieuwe bestelling"
$logs.Range("B$newRow").Value = "planning@testbedrijf123.nl"
$logs.Range("C$newRow").Value = "Wil je 200 stuks M8-bouten bestellen bij onze leverancier?"
$logs.Range("D$newRow").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E$newRow").Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@testbedrijf123.nl."
$logs.Range("F$newRow").Value = "2025-08-14 21:57:54"
$logs.Range("G$newRow").Value = "Nee"
$logs.Range("H$newRow").Value = "Ja"
$logs.Range("I$newRow").Value = "Nee"
$logs.Range("J$newRow").Value = "Nee"

# Extend the conditional formatting ranges (D, G, H, I, J) from row 33 to row 34
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "33")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "$newRow")
    $conditions = $oldRange.FormatConditions
    if ($conditions.Count -gt 0) {
        $conditions.Item(1).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard summary count for "Intern verzoek / Actie voor medewerker"
$dashboard.Range("B2").Value = 26
